# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型"
# sheets for a handful of rows, reflecting the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 13902
    $ws.Range("F11").Value = 65
    $ws.Range("F15").Value = 5910
    $ws.Range("F17").Value = 94
}
